$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1): "_old" -> "_FV2404", "_new" -> "_FV2410" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$oldNames = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $oldNames[$i]
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$newNames = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $newNames[$i]
}

# --- Turn the used range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U76"), 0, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
